# Update the MSME Country Indicators - Mozambique Summary figures.
# The target cells hold their numbers as plain *text* (shared strings),
# not numeric values, so we temporarily force a text number-format before
# assigning the new value (otherwise Excel auto-converts "1.13" etc. into
# a real number) and then restore the cell's original style/appearance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewText
    )

    $range = $ws.Range($CellRef)
    $originalStyle = $range.Style

    $range.NumberFormat = "@"
    $range.Value = $NewText
    $range.Style = $originalStyle
}

# Enterprises density (per 1000 people): row 13
Set-TextValue "B13" "1.13"
Set-TextValue "C13" "0.26"
Set-TextValue "D13" "1.39"

# Employment (% of total): row 14
Set-TextValue "B14" "14.13"
Set-TextValue "C14" "28.78"
Set-TextValue "D14" "42.91"

# Enterprises (% of total): row 16
Set-TextValue "B16" "71.13"
Set-TextValue "C16" "16.65"
Set-TextValue "D16" "87.77"
